$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D receive numeric-looking text; force text storage
# by switching the cell to Text format before assignment, then restore
# the default "Normal" style so no residual style index is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.170.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.886.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "490.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.01%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.415"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.09"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.344"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.377.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.45%  "

$ws.Range("E13").Value = "  -4.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.79%  "

$ws.Range("E15").Value = "  -7.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "55.125.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.73%  "

$ws.Range("E17").Value = "  -5.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.886.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.476"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "61.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.002.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.157"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0827"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.37%  "

$ws.Range("E34").Value = "  -10.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "148.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.25%  "

$ws.Range("E39").Value = "  -10.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0642"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.627"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.070.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.902"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.15%  "

$ws.Range("E49").Value = "  -5.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0831"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.52%  "
